# Daily attendance processing - re-sort the "Recorded By" (column G) list of
# recorders for every data row so the comma-separated names appear in
# ordinal (case-sensitive, ASCII) alphabetical order, e.g.
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#
# Note: Sort-Object in this host does case-INsensitive comparisons, which
# does not reproduce the required ordinal ordering (uppercase before
# lowercase). So we do a small in-place insertion sort using .CompareTo(),
# which behaves ordinally here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = @($val -split ", ")
    $n = $parts.Count

    # Insertion sort (ordinal / case-sensitive via .CompareTo) - lists are
    # short (2-3 entries) so this is plenty fast.
    for ($i = 1; $i -lt $n; $i++) {
        $key = $parts[$i]
        $j = $i - 1
        while ($j -ge 0 -and $parts[$j].CompareTo($key) -gt 0) {
            $parts[$j + 1] = $parts[$j]
            $j = $j - 1
        }
        $parts[$j + 1] = $key
    }

    $newVal = $parts -join ", "
    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
